$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Mention" values (col L) for rows 26-101, replacing the placeholder
# "À recalculer" with the recalculated value (mirrors the source commit
# "Recalcul automatique de la colonne Mention pour les anciennes interventions").
$mentions = @(
    "Élu",
    "Élu",
    "Élu",
    "Élu",
    "Élu",
    "Élu & Conseil fédéral",
    "Conseil fédéral",
    "Élu",
    "Conseil fédéral",
    "Conseil fédéral",
    "Conseil fédéral",
    "Élu",
    "Élu & Conseil fédéral",
    "Élu",
    "Élu",
    "Élu & Conseil fédéral",
    "Conseil fédéral",
    "Conseil fédéral",
    "Conseil fédéral",
    "Conseil fédéral",
    "Élu & Conseil fédéral",
    "Conseil fédéral",
    "Conseil fédéral",
    "Élu",
    "Élu",
    "Élu",
    "Conseil fédéral",
    "Conseil fédéral",
    "Élu",
    "Conseil fédéral",
    "Conseil fédéral",
    "Élu & Conseil fédéral",
    "Conseil fédéral",
    "Conseil fédéral",
    "Conseil fédéral",
    "Élu",
    "Élu",
    "Conseil fédéral",
    "Élu",
    "Conseil fédéral",
    "Élu & Conseil fédéral",
    "Conseil fédéral",
    "Élu & Conseil fédéral",
    "Élu & Conseil fédéral",
    "Élu & Conseil fédéral",
    "Élu",
    "Conseil fédéral",
    "Élu",
    "Conseil fédéral",
    "Élu & Conseil fédéral",
    "Élu & Conseil fédéral",
    "Élu",
    "Conseil fédéral",
    "Conseil fédéral",
    "Conseil fédéral",
    "Élu & Conseil fédéral",
    "Conseil fédéral",
    "Élu",
    "Conseil fédéral",
    "Élu",
    "Conseil fédéral",
    "Élu",
    "Conseil fédéral",
    "Conseil fédéral",
    "Élu & Conseil fédéral",
    "Conseil fédéral",
    "Conseil fédéral",
    "Conseil fédéral",
    "Conseil fédéral",
    "Élu & Conseil fédéral",
    "Élu",
    "Élu",
    "Conseil fédéral",
    "Élu",
    "Conseil fédéral",
    "Élu"
)

$startRow = 26
for ($i = 0; $i -lt $mentions.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 12).Value = $mentions[$i]
}
